# Update "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the 6443f018-... row on the zh-cn and de-de report sheets,
# reflecting the refreshed handback report generation timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-23 08:49:26"
$wsZhCn.Range("H4").Value = "2016-03-23 08:49:51"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-23 08:49:30"
$wsDeDe.Range("H4").Value = "2016-03-23 08:49:58"
